$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "email" value cells for John and Jack, linked as mailto hyperlinks.
# Setting the cell text first keeps the visible/shared-string text as the
# plain email address, and Hyperlinks.Add then wires up the mailto: target
# plus the built-in "Hyperlink" style (new font + cellStyleXfs/cellXfs/
# cellStyles entries), matching the OOXML diff.
$ws.Range("G3").Value = "john@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:john@gmail.com")

$ws.Range("G4").Value = "jack@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:jack@gmail.com")

# Update the active selection to J5, matching the recorded selection change.
[void]$ws.Range("J5").Select()
